$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Final roster table (row, Player Name, Position, Team)
$data = @(
    @(2,  "Jalen Brunson",      "PG",       "New York Knicks"),
    @(3,  "Coby White",         "PG,SG",    "Chicago Bulls"),
    @(4,  "Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @(5,  "Shaedon Sharpe",     "SG,SF",    "Portland Trail Blazers"),
    @(6,  "LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @(7,  "Desmond Bane",       "SG,SF",    "Memphis Grizzlies"),
    @(8,  "Dalton Knecht",      "SG,SF",    "Los Angeles Lakers"),
    @(9,  "Nicolas Claxton",    "C",        "Brooklyn Nets"),
    @(10, "Dereck Lively II",   "C",        "Dallas Mavericks"),
    @(11, "Moussa Diabate",     "C",        "Charlotte Hornets"),
    @(12, "Trae Young",         "PG",       "Atlanta Hawks"),
    @(13, "Devin Booker",       "PG,SG",    "Phoenix Suns"),
    @(14, "P.J. Washington",    "PF",       "Dallas Mavericks"),
    @(15, "Alperen Sengün",     "C",        "Houston Rockets"),
    @(16, "Walker Kessler",     "C",        "Utah Jazz"),
    @(17, "Immanuel Quickley",  "PG,SG",    "Toronto Raptors"),
    @(18, "Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @(19, "Norman Powell",      "SG,SF",    "LA Clippers")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
